$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" summary text with the day's updated conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.61 = 5935.1 pesos`n✅ 5935.1 pesos = 1.6 = 962.97 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the "tasas" sheet rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("O10").Value = 3679.76
$ws2.Range("N12").Value = 3698
$ws2.Range("O12").Value = 600.001
